$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Version" column only within the property-definitions sub-table
# (rows 4-8, column A), shifting the remaining cells in those rows to the left.
$ws.Range("A4:A8").Delete(-4159)  # -4159 = xlShiftToLeft

# Select A4 to match final selection state
$ws.Range("A4").Select()
